# Updated cryptos list on Mon Mar 27 10:15:45 UTC 2023 with GitHub Actions
#
# The "Price" column (D) and "Volume(1h)" column (E) hold plain text in the
# sheet (coinranking.com's formatted strings, e.g. "27.948.32" or
# "  +0.79%  "), not real numbers. Several of the new Price strings
# (e.g. "328.76", "1.000", "0.4680", "20.60") look numeric, so a naive
# Range.Value assignment would make Excel auto-convert them to actual
# numbers (dropping the literal formatting and the original cell style).
# Forcing the cell to Text format before the write, then restoring the
# "Normal" style afterwards, keeps the write a literal text value while
# leaving the cell's style index untouched (matches the original, unstyled
# cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $text) {
    $cell = $ws.Range($a1)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$updates = [ordered]@{
    2  = @{ D = "27.948.32";  E = "  +0.79%  " }
    3  = @{ D = "1.764.75";   E = "  -0.64%  " }
    4  = @{                   E = "  -0.05%  " }
    5  = @{ D = "328.76";     E = "  +0.51%  " }
    6  = @{ D = "1.000";      E = "  -0.07%  " }
    7  = @{ D = "0.4680";     E = "  +1.30%  " }
    8  = @{ D = "0.3516";     E = "  -1.84%  " }
    9  = @{ D = "42.93";      E = "  +2.59%  " }
    10 = @{ D = "0.07366";    E = "  -1.56%  " }
    11 = @{ D = "1.081";      E = "  -1.89%  " }
    12 = @{ D = "1.001";      E = "  +0.00%  " }
    13 = @{ D = "20.60";      E = "  -1.05%  " }
    14 = @{ D = "5.999";      E = "  -0.69%  " }
    15 = @{ D = "7.167";      E = "  -0.82%  " }
    16 = @{ D = "1.762.09";   E = "  -0.85%  " }
    17 = @{ D = "92.35" }
    18 = @{ D = "0.00001054"; E = "  -0.39%  " }
    19 = @{ D = "0.06422";    E = "  +0.19%  " }
    20 = @{ D = "0.9999";     E = "  -0.06%  " }
    21 = @{ D = "16.87";      E = "  -1.26%  " }
    22 = @{ D = "5.762";      E = "  -0.36%  " }
    23 = @{ D = "27.975.91";  E = "  +0.62%  " }
    24 = @{                   E = "  -1.30%  " }
    25 = @{ D = "2.150";      E = "  +3.14%  " }
    26 = @{ D = "162.15";     E = "  -1.24%  " }
    27 = @{                   E = "  -1.18%  " }
    28 = @{ D = "1.965.03";   E = "  -0.75%  " }
    29 = @{ D = "2.180";      E = "  +0.58%  " }
    30 = @{ D = "123.06";     E = "  -2.26%  " }
    31 = @{ D = "1.070";      E = "  -1.88%  " }
    32 = @{ D = "0.09289";    E = "  +0.59%  " }
    33 = @{ D = "3.646";      E = "  -0.65%  " }
    34 = @{ D = "5.549";      E = "  +0.34%  " }
    35 = @{ D = "11.68";      E = "  -1.09%  " }
    36 = @{ D = "0.02265";    E = "  -1.21%  " }
    37 = @{ D = "0.06069";    E = "  -1.69%  " }
    38 = @{ D = "0.2060";     E = "  -1.15%  " }
    39 = @{ D = "4.910";      E = "  -0.79%  " }
    40 = @{ D = "0.6123";     E = "  -2.91%  " }
    41 = @{ D = "1.184";      E = "  -0.09%  " }
    42 = @{ D = "1.387";      E = "  -0.31%  " }
    43 = @{ D = "7.759";      E = "  -0.22%  " }
    44 = @{ D = "13.13";      E = "  -1.07%  " }
    45 = @{ D = "3.740";      E = "  +0.06%  " }
    46 = @{ D = "0.5783";     E = "  -1.68%  " }
    47 = @{ D = "123.22";     E = "  +0.78%  " }
    48 = @{ D = "1.926";      E = "  -1.01%  " }
    49 = @{ D = "0.06823";    E = "  -1.64%  " }
    50 = @{ D = "1.123";      E = "  -1.12%  " }
    51 = @{ D = "72.00";      E = "  -0.18%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        Set-TextValue "$col$row" $cols[$col]
    }
}
